$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first block of changes
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14294
$ws1.Range("F3").Value = 336
$ws1.Range("F6").Value = 563
$ws1.Range("F7").Value = 1493
$ws1.Range("F8").Value = 143

# Sheet "全部类型" (All types) - second block of changes
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14294
$ws4.Range("F3").Value = 336
$ws4.Range("F8").Value = 563
$ws4.Range("F9").Value = 1493
$ws4.Range("F11").Value = 143
